$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, reusing the same formatting (bold font + border + centered
# alignment) as the rest of the header row by copying H1's format onto I1.
$ws.Range("I1").Value = "Rate of exciton generation"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Data cells: Rate of exciton generation = (# of photons) / 50
$ws.Range("I2").Value = 3703965503269865000000.0
$ws.Range("I3").Value = 3710990277995830000000.0
$ws.Range("I4").Value = 3718595265851428000000.0
$ws.Range("I5").Value = 3726776172055099000000.0
$ws.Range("I6").Value = 3735528374983967000000.0
$ws.Range("I7").Value = 3744846928496203000000.0
$ws.Range("I8").Value = 3754726564442176000000.0
